$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').Value = '63.987.23'
$ws.Range('E2').Value = '  -5.37%  '
$ws.Range('D3').Value = '3.289.13'
$ws.Range('E3').Value = '  -6.61%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '522.48'
$ws.Range('E5').Value = '  -5.30%  '
$ws.Range('D6').Value = '172.42'
$ws.Range('E6').Value = '  -14.51%  '
$ws.Range('D7').Value = '0.600'
$ws.Range('E7').Value = '  -1.92%  '
$ws.Range('D8').Value = '3.286.77'
$ws.Range('E8').Value = '  -6.39%  '
$ws.Range('E9').Value = '  +0.06%  '
$ws.Range('D10').Value = '0.604'
$ws.Range('E10').Value = '  -7.51%  '
$ws.Range('D11').Value = '55.81'
$ws.Range('E11').Value = '  -11.10%  '
$ws.Range('D12').Value = '0.132'
$ws.Range('E12').Value = '  -7.66%  '
$ws.Range('E13').Value = '  -5.03%  '
$ws.Range('E14').Value = '  -8.40%  '
$ws.Range('D15').Value = '3.818.14'
$ws.Range('E15').Value = '  -6.62%  '
$ws.Range('D16').Value = '3.290.68'
$ws.Range('E16').Value = '  -6.67%  '
$ws.Range('D17').Value = '0.115'
$ws.Range('E17').Value = '  -6.97%  '
$ws.Range('D18').Value = '63.938.10'
$ws.Range('E18').Value = '  -5.16%  '
$ws.Range('D19').Value = '17.34'
$ws.Range('E19').Value = '  -5.97%  '
$ws.Range('D20').Value = '11.02'
$ws.Range('E20').Value = '  -6.89%  '
$ws.Range('D21').Value = '0.957'
$ws.Range('E21').Value = '  -6.72%  '
$ws.Range('D22').Value = '371.93'
$ws.Range('E22').Value = '  -5.46%  '
$ws.Range('E23').Value = '  -5.98%  '
$ws.Range('D24').Value = '80.12'
$ws.Range('E24').Value = '  -4.72%  '
$ws.Range('D25').Value = '10.93'
$ws.Range('E25').Value = '  -8.72%  '
$ws.Range('D26').Value = '3.89'
$ws.Range('E26').Value = '  -0.05%  '
$ws.Range('B27').Value = 'ImmutableX'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D27').Value = '2.68'
$ws.Range('E27').Value = '  -5.30%  '
$ws.Range('B28').Value = 'InternetComputer(DFINITY)'
$ws.Range('C28').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D28').Value = '11.29'
$ws.Range('E28').Value = '  -7.84%  '
$ws.Range('B29').Value = 'Filecoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D29').Value = '8.27'
$ws.Range('E29').Value = '  -6.57%  '
$ws.Range('B30').Value = 'EthereumClassic'
$ws.Range('C30').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D30').Value = '28.66'
$ws.Range('E30').Value = '  -7.86%  '
$ws.Range('B31').Value = 'Bittensor'
$ws.Range('C31').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D31').Value = '639.41'
$ws.Range('E31').Value = '  -11.09%  '
$ws.Range('B32').Value = 'NEARProtocol'
$ws.Range('C32').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D32').Value = '6.58'
$ws.Range('E32').Value = '  -7.00%  '
$ws.Range('B33').Value = 'Cosmos'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D33').Value = '11.19'
$ws.Range('E33').Value = '  -4.67%  '
$ws.Range('B34').Value = 'OKB'
$ws.Range('C34').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D34').Value = '58.81'
$ws.Range('E34').Value = '  -7.57%  '
$ws.Range('E35').Value = '  -5.95%  '
$ws.Range('B36').Value = 'Dai'
$ws.Range('C36').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D36').Value = '1.00'
$ws.Range('E36').Value = '  -0.02%  '
$ws.Range('B37').Value = 'InjectiveProtocol'
$ws.Range('C37').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D37').Value = '36.67'
$ws.Range('E37').Value = '  -4.82%  '
$ws.Range('B38').Value = 'TheGraph'
$ws.Range('C38').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D38').Value = '0.383'
$ws.Range('E38').Value = '  -3.43%  '
$ws.Range('B39').Value = 'FirstDigitalUSD'
$ws.Range('C39').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D39').Value = '0.998'
$ws.Range('E39').Value = '  -0.04%  '
$ws.Range('B40').Value = 'PEPE'
$ws.Range('C40').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D40').Value = '0.0₃0691'
$ws.Range('E40').Value = '  +1.65%  '
$ws.Range('B41').Value = 'Maker'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D41').Value = '2.923.17'
$ws.Range('E41').Value = '  -5.13%  '
$ws.Range('B42').Value = 'Kaspa'
$ws.Range('C42').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D42').Value = '0.121'
$ws.Range('E42').Value = '  -7.65%  '
$ws.Range('B43').Value = 'Fetch.AI'
$ws.Range('C43').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D43').Value = '2.43'
$ws.Range('E43').Value = '  -4.40%  '
$ws.Range('B44').Value = 'ThetaToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D44').Value = '2.67'
$ws.Range('E44').Value = '  -12.17%  '
$ws.Range('B45').Value = 'WEMIXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D45').Value = '2.63'
$ws.Range('E45').Value = '  -6.44%  '
$ws.Range('D46').Value = '0.0394'
$ws.Range('E46').Value = '  -3.63%  '
$ws.Range('B47').Value = 'ApeXProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D47').Value = '2.99'
$ws.Range('E47').Value = '  +3.71%  '
$ws.Range('B48').Value = 'Stellar'
$ws.Range('C48').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D48').Value = '0.124'
$ws.Range('E48').Value = '  -2.47%  '
$ws.Range('B49').Value = 'Stacks'
$ws.Range('C49').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D49').Value = '2.73'
$ws.Range('E49').Value = '  +3.89%  '
$ws.Range('B50').Value = 'Monero'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D50').Value = '134.91'
$ws.Range('E50').Value = '  -2.40%  '
$ws.Range('B51').Value = 'dogwifhat'
$ws.Range('C51').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D51').Value = '2.35'
$ws.Range('E51').Value = '  -12.61%  '
